$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing CombinedQL row (row 4) values
$ws.Range("B4").Value = 95
$ws.Range("C4").Value = 0.6251507870414386
$ws.Range("D4").Value = 0.01806664718814597
$ws.Range("E4").Value = 0.845321769902543

# Update DualQL row (row 5): only D5 changes
$ws.Range("D5").Value = 0.006999774671284044

# Add new row 6 for DWA, matching the formatting used by the other label cells
$ws.Range("A6").Value = "DWA"
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$ws.Range("B6").Value = 80
$ws.Range("C6").Value = 0.7573459422768216
$ws.Range("D6").Value = 0.02719741163277957
$ws.Range("E6").Value = 0.7429682400799604
